$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("KN1").Value = 44136
$ws.Range("KO1").Value = 44137
$ws.Range("KP1").Value = 44138
$ws.Range("KQ1").Value = 44139
$ws.Range("KR1").Value = 44140
$ws.Range("KS1").Value = 44141
$ws.Range("KT1").Value = 44142
$ws.Range("KU1").Value = 44143
$ws.Range("KV1").Value = 44144
$ws.Range("KW1").Value = 44145

$ws.Range("KN2").Value = 43.1
$ws.Range("KO2").Value = 47.29
$ws.Range("KP2").Value = 43.62
$ws.Range("KQ2").Value = 44.2
$ws.Range("KR2").Value = 50.18
$ws.Range("KS2").Value = 58.54
$ws.Range("KT2").Value = 62.59
$ws.Range("KU2").Value = 50.02
$ws.Range("KV2").Value = 49.65
$ws.Range("KW2").Value = 51.29

$ws.Range("KN3").Value = 33.06
$ws.Range("KO3").Value = 35.01
$ws.Range("KP3").Value = 34.32
$ws.Range("KQ3").Value = 34.33
$ws.Range("KR3").Value = 35.78
$ws.Range("KS3").Value = 37.87
$ws.Range("KT3").Value = 34.67
$ws.Range("KU3").Value = 31.31
$ws.Range("KV3").Value = 34.04
$ws.Range("KW3").Value = 32.43

$ws.Range("KN4").Value = 48.78
$ws.Range("KO4").Value = 64.13
$ws.Range("KP4").Value = 58.52
$ws.Range("KQ4").Value = 58.93
$ws.Range("KR4").Value = 59.67
$ws.Range("KS4").Value = 61.96
$ws.Range("KT4").Value = 60.73
$ws.Range("KU4").Value = 52.66
$ws.Range("KV4").Value = 67.3
$ws.Range("KW4").Value = 60.78

$ws.Range("KN5").Value = 52.72
$ws.Range("KO5").Value = 62.59
$ws.Range("KP5").Value = 61.18
$ws.Range("KQ5").Value = 59.31
$ws.Range("KR5").Value = 59.72
$ws.Range("KS5").Value = 67.99
$ws.Range("KT5").Value = 60.38
$ws.Range("KU5").Value = 49.36
$ws.Range("KV5").Value = 60.69
$ws.Range("KW5").Value = 61.92

$ws.Range("KM1").Copy() | Out-Null
$ws.Range("KN1:KW1").PasteSpecial(-4122) | Out-Null

$ws.Range("KW1").Select() | Out-Null
